$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1600.9474
$ws.Range("I15").Value = 1600.9474
$ws.Range("K15").Value = 4802.8422
$ws.Range("M15").Value = -4633.8422

$ws.Range("H87").Value = 20627
$ws.Range("J87").Value = 20627
$ws.Range("L87").Value = 20627
$ws.Range("N87").Value = -23123

$ws.Range("H88").Value = 6416.5
$ws.Range("I88").Value = 7750
$ws.Range("J88").Value = 3749.5
$ws.Range("K88").Value = 7750
$ws.Range("L88").Value = 3749.5
$ws.Range("M88").Value = -7344
$ws.Range("N88").Value = -4561.5

$ws.Range("H90").Value = 20627
$ws.Range("J90").Value = 20627
$ws.Range("L90").Value = 61881
$ws.Range("N90").Value = -74361

$ws.Range("H91").Value = 6416.5
$ws.Range("I91").Value = 7750
$ws.Range("J91").Value = 3749.5
$ws.Range("K91").Value = 7750
$ws.Range("L91").Value = 3749.5
$ws.Range("M91").Value = -6346
$ws.Range("N91").Value = -6557.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 5000000
$ws.Range("I8").Value = 5000000
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 5000000
$ws.Range("L8").Value = 0
$ws.Range("N8").Value = -4999856
$ws.Range("M8").ClearContents()

$ws.Range("H32").Value = 11769554
$ws.Range("I32").Value = 12200047
$ws.Range("K32").Value = 12200047
$ws.Range("M32").Value = -12199760

$ws.Range("H45").Value = 4019.6155
$ws.Range("I45").Value = 3045.5557
$ws.Range("J45").Value = 6211.25
$ws.Range("K45").Value = 3045.5557
$ws.Range("L45").Value = 6211.25
$ws.Range("M45").Value = -2668.5557
$ws.Range("N45").Value = -6965.25

$ws.Range("H102").Value = 2188.125
$ws.Range("I102").Value = 2329.5833
$ws.Range("J102").Value = 1763.75
$ws.Range("K102").Value = 2329.5833
$ws.Range("L102").Value = 1763.75
$ws.Range("M102").Value = -707.5832999999998
$ws.Range("N102").Value = -5007.75

$ws.Range("H110").Value = 992.2857
$ws.Range("I110").Value = 992.2857
$ws.Range("K110").Value = 992.2857
$ws.Range("M110").Value = 1052.7143

$ws.Range("H122").Value = 3759.6086
$ws.Range("I122").Value = 2863.8667
$ws.Range("J122").Value = 5439.125
$ws.Range("K122").Value = 8591.6001
$ws.Range("L122").Value = 16317.375
$ws.Range("M122").Value = -6141.6001
$ws.Range("N122").Value = -21217.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1896.1111
$ws.Range("I99").Value = 1826
$ws.Range("J99").Value = 2015.3
$ws.Range("K99").Value = 1826
$ws.Range("L99").Value = 2015.3
$ws.Range("M99").Value = -328
$ws.Range("N99").Value = -5011.3

$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("N107").Value = 0
$ws.Range("L107").ClearContents()
$ws.Range("M107").ClearContents()

$ws.Range("H134").Value = 3670.5227
$ws.Range("I134").Value = 2526.6785
$ws.Range("J134").Value = 5672.25
$ws.Range("K134").Value = 7580.0355
$ws.Range("L134").Value = 17016.75
$ws.Range("M134").Value = -5045.0355
$ws.Range("N134").Value = -22086.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 10000
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 10000
$ws.Range("K11").Value = 0
$ws.Range("M11").Value = 10000
$ws.Range("N11").Value = -10280
$ws.Range("L11").ClearContents()

$ws.Range("H31").Value = 3520
$ws.Range("I31").Value = 2496.5293
$ws.Range("J31").Value = 6999.8
$ws.Range("K31").Value = 2496.5293
$ws.Range("L31").Value = 6999.8
$ws.Range("M31").Value = -2201.5293
$ws.Range("N31").Value = -7589.8

$ws.Range("H34").Value = 3520
$ws.Range("I34").Value = 2496.5293
$ws.Range("J34").Value = 6999.8
$ws.Range("K34").Value = 2496.5293
$ws.Range("L34").Value = 6999.8
$ws.Range("M34").Value = -2294.5293
$ws.Range("N34").Value = -7403.8

$ws.Range("H62").Value = 8499.333000000001
$ws.Range("J62").Value = 3166.3333
$ws.Range("L62").Value = 3166.3333
$ws.Range("N62").Value = -4414.3333

$ws.Range("H65").Value = 8499.333000000001
$ws.Range("J65").Value = 3166.3333
$ws.Range("L65").Value = 15831.6665
$ws.Range("N65").Value = -22071.6665

$ws.Range("H134").Value = 3514.0298
$ws.Range("I134").Value = 2299.6445
$ws.Range("J134").Value = 5998
$ws.Range("K134").Value = 6898.933499999999
$ws.Range("L134").Value = 17994
$ws.Range("M134").Value = -4363.933499999999
$ws.Range("N134").Value = -23064

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 954
$ws.Range("I7").Value = 1128.1333
$ws.Range("J7").Value = 301
$ws.Range("K7").Value = 3384.3999
$ws.Range("L7").Value = 903
$ws.Range("M7").Value = -3272.3999
$ws.Range("N7").Value = -1127

$ws.Range("H18").Value = 858.8333
$ws.Range("I18").Value = 1030
$ws.Range("J18").Value = 824.6
$ws.Range("K18").Value = 3090
$ws.Range("L18").Value = 2473.8
$ws.Range("M18").Value = -2921
$ws.Range("N18").Value = -2811.8

$ws.Range("H61").Value = 294.76923
$ws.Range("I61").Value = 183.2
$ws.Range("J61").Value = 666.6667
$ws.Range("K61").Value = 549.5999999999999
$ws.Range("L61").Value = 2000.0001
$ws.Range("M61").Value = -334.5999999999999
$ws.Range("N61").Value = -2430.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 661.1429000000001
$ws.Range("I13").Value = 809.6667
$ws.Range("J13").Value = 549.75
$ws.Range("K13").Value = 809.6667
$ws.Range("L13").Value = 549.75
$ws.Range("M13").Value = -670.6667
$ws.Range("N13").Value = -827.75

$ws.Range("H97").Value = 947.4
$ws.Range("I97").Value = 750
$ws.Range("J97").Value = 996.75
$ws.Range("K97").Value = 750
$ws.Range("L97").Value = 996.75
$ws.Range("M97").Value = -254
$ws.Range("N97").Value = -1988.75

$ws.Range("H102").Value = 37905.184
$ws.Range("I102").Value = 42388.758
$ws.Range("J102").Value = 5399.25
$ws.Range("K102").Value = 42388.758
$ws.Range("L102").Value = 5399.25
$ws.Range("M102").Value = -40766.758
$ws.Range("N102").Value = -8643.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 42448.145
$ws.Range("I100").Value = 124525.11
$ws.Range("J100").Value = 3569.5789
$ws.Range("K100").Value = 124525.11
$ws.Range("L100").Value = 3569.5789
$ws.Range("M100").Value = -123984.11
$ws.Range("N100").Value = -4651.5789

$ws.Range("H139").Value = 84411
$ws.Range("J139").Value = 83437.19
$ws.Range("L139").Value = 83437.19
$ws.Range("N139").Value = -93717.19

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1736.4324
$ws.Range("I132").Value = 1113.8518
$ws.Range("J132").Value = 3417.4
$ws.Range("K132").Value = 3341.5554
$ws.Range("L132").Value = 10252.2
$ws.Range("M132").Value = -811.5553999999997
$ws.Range("N132").Value = -15312.2

